# ordersList.xlsx - rename three customers (得意先名, column D) and widen
# the column to fit the new (longer) names, then park the selection on D1.
#
#   伊藤商事(Ito Shoji)       -> さとう商事(Sato Shoji)
#   ライトオフ(Light Off)     -> ノーズライト(Noselight)
#   ビックマックハウス(Big Mac House) -> ビックハウス(Big House)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (1-based) in column D currently holding each customer name, taken
# from the "得意先名" column of the orders list. Writing in this order
# reproduces the shared-string insertion order of the source edit.
$noseLightRows = @(6, 7, 8, 9, 10, 20, 21, 22, 23, 24, 28, 29, 30, 31)
$bigHouseRows  = @(13, 14, 15, 16, 17, 18)
$satoShojiRows = @(2, 3, 4, 5)

foreach ($r in $noseLightRows) { $ws.Cells($r, 4).Value = "ノーズライト" }
foreach ($r in $bigHouseRows)  { $ws.Cells($r, 4).Value = "ビックハウス" }
foreach ($r in $satoShojiRows) { $ws.Cells($r, 4).Value = "さとう商事" }

# The new names are longer than "OSAKA BASE" / "ビックマックハウス", so
# widen column D to fit them.
$ws.Columns(4).ColumnWidth = 28.5

# Park the view on D1 (previously scrolled to row 11 with E37 selected).
$ws.Range("D1").Select()
